# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de):
#  - Status (col C) moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Latest Target File (col F) gets the same source .md file as col A (with hyperlink)
#  - Latest Handback File (col G) gets the same target .xlf file as col D (with hyperlink)
#  - Latest Handback DateTime (col H) gets stamped with the handback time
# The Overview sheet's Status columns (B/C) mirror the same text.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: reflect the new status text for both rows ---
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Range("B2").Value = $statusText
$wsOv.Range("C2").Value = $statusText
$wsOv.Range("B3").Value = $statusText
$wsOv.Range("C3").Value = $statusText

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Hyperlinks.Add(
    $wsZh.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/c35cd6b5923e679c4f7412205471f52430c2663f/e2e/03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.md",
    "",
    "",
    "03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/920ef8d2b562787dd1fee2b024f490a10af13367/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.58559e9ba60de1b8adcfaef4d30383980f106c5c.zh-cn.xlf",
    "",
    "",
    "03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.58559e9ba60de1b8adcfaef4d30383980f106c5c.zh-cn.xlf"
) | Out-Null
$wsZh.Range("H2").Value = "2016-03-24 16:58:23"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/c35cd6b5923e679c4f7412205471f52430c2663f/e2e/2ef9fcc0-8841-4f9c-85f6-6390a597d252.md",
    "",
    "",
    "2ef9fcc0-8841-4f9c-85f6-6390a597d252.md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/920ef8d2b562787dd1fee2b024f490a10af13367/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2ef9fcc0-8841-4f9c-85f6-6390a597d252.469605ca99d80052814dd73d9ea6d6c62f405d18.zh-cn.xlf",
    "",
    "",
    "2ef9fcc0-8841-4f9c-85f6-6390a597d252.469605ca99d80052814dd73d9ea6d6c62f405d18.zh-cn.xlf"
) | Out-Null
$wsZh.Range("H3").Value = "2016-03-24 16:58:23"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Hyperlinks.Add(
    $wsDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/c35cd6b5923e679c4f7412205471f52430c2663f/e2e/03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.md",
    "",
    "",
    "03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cf8c1f1c5b476c5c3fbdf470cb8d78f54478e0ac/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.58559e9ba60de1b8adcfaef4d30383980f106c5c.de-de.xlf",
    "",
    "",
    "03f9a1cb-e8aa-4130-baf0-57c4d76a84e1.58559e9ba60de1b8adcfaef4d30383980f106c5c.de-de.xlf"
) | Out-Null
$wsDe.Range("H2").Value = "2016-03-24 16:58:34"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/c35cd6b5923e679c4f7412205471f52430c2663f/e2e/2ef9fcc0-8841-4f9c-85f6-6390a597d252.md",
    "",
    "",
    "2ef9fcc0-8841-4f9c-85f6-6390a597d252.md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cf8c1f1c5b476c5c3fbdf470cb8d78f54478e0ac/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2ef9fcc0-8841-4f9c-85f6-6390a597d252.469605ca99d80052814dd73d9ea6d6c62f405d18.de-de.xlf",
    "",
    "",
    "2ef9fcc0-8841-4f9c-85f6-6390a597d252.469605ca99d80052814dd73d9ea6d6c62f405d18.de-de.xlf"
) | Out-Null
$wsDe.Range("H3").Value = "2016-03-24 16:58:34"

Write-Host "Handback report generated"
